$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# The "Final Amount" header row (row 1) and its result column (E) are no
# longer needed - remove the header row and drop the now-unused column E.
$ws1.Rows.Item(1).Delete()
$ws1.Columns.Item(5).Delete()

# Excel recalculates the default row height for the remaining data rows
# once the taller/bolder header row is gone.
$ws1.Range("A1:D5").RowHeight = 17.9
